$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) and Volume(1h) (column E) updates for the cryptos list refresh.
# Column D values are written via a Text number format so Excel does not
# reinterpret numeric-looking strings (e.g. "0.999") as numbers, matching the
# original inline-string cell content; the format is then reset to Normal so
# no visible style change is left behind.

$d = $ws.Cells.Item(2, 4)
$d.NumberFormat = "@"
$d.Value = "60.584.85"
$d.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  -2.36%  "

$d = $ws.Cells.Item(3, 4)
$d.NumberFormat = "@"
$d.Value = "2.903.21"
$d.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  -3.26%  "

$ws.Cells.Item(4, 5).Value = "  +0.07%  "

$d = $ws.Cells.Item(5, 4)
$d.NumberFormat = "@"
$d.Value = "527.71"
$d.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -4.18%  "

$d = $ws.Cells.Item(6, 4)
$d.NumberFormat = "@"
$d.Value = "142.93"
$d.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -7.66%  "

$d = $ws.Cells.Item(7, 4)
$d.NumberFormat = "@"
$d.Value = "0.999"
$d.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  -0.07%  "

$d = $ws.Cells.Item(8, 4)
$d.NumberFormat = "@"
$d.Value = "0.556"
$d.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -1.78%  "

$d = $ws.Cells.Item(9, 4)
$d.NumberFormat = "@"
$d.Value = "2.911.37"
$d.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -3.17%  "

$ws.Cells.Item(10, 5).Value = "  -3.89%  "

$ws.Cells.Item(11, 5).Value = "  -5.34%  "

$d = $ws.Cells.Item(12, 4)
$d.NumberFormat = "@"
$d.Value = "0.360"
$d.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -1.90%  "

$d = $ws.Cells.Item(13, 4)
$d.NumberFormat = "@"
$d.Value = "3.406.62"
$d.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -3.37%  "

$ws.Cells.Item(14, 5).Value = "  +1.54%  "

$d = $ws.Cells.Item(15, 4)
$d.NumberFormat = "@"
$d.Value = "60.571.82"
$d.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -2.47%  "

$d = $ws.Cells.Item(16, 4)
$d.NumberFormat = "@"
$d.Value = "22.58"
$d.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -5.01%  "

$d = $ws.Cells.Item(17, 4)
$d.NumberFormat = "@"
$d.Value = "2.905.97"
$d.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -3.15%  "

$d = $ws.Cells.Item(18, 4)
$d.NumberFormat = "@"
$d.Value = "0.0000142"
$d.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  -4.50%  "

$d = $ws.Cells.Item(19, 4)
$d.NumberFormat = "@"
$d.Value = "5.03"
$d.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -1.78%  "

$ws.Cells.Item(20, 5).Value = "  -2.83%  "

$d = $ws.Cells.Item(21, 4)
$d.NumberFormat = "@"
$d.Value = "363.62"
$d.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -7.73%  "

$d = $ws.Cells.Item(22, 4)
$d.NumberFormat = "@"
$d.Value = "6.57"
$d.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -1.16%  "

$ws.Cells.Item(23, 5).Value = "  -0.09%  "

$d = $ws.Cells.Item(24, 4)
$d.NumberFormat = "@"
$d.Value = "63.61"
$d.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -2.41%  "

$d = $ws.Cells.Item(25, 4)
$d.NumberFormat = "@"
$d.Value = "3.018.75"
$d.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -3.33%  "

$ws.Cells.Item(26, 5).Value = "  -3.66%  "

$d = $ws.Cells.Item(27, 4)
$d.NumberFormat = "@"
$d.Value = "0.181"
$d.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -2.21%  "

$d = $ws.Cells.Item(28, 4)
$d.NumberFormat = "@"
$d.Value = "0.998"
$d.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -0.02%  "

$d = $ws.Cells.Item(29, 4)
$d.NumberFormat = "@"
$d.Value = "7.84"
$d.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -7.47%  "

$d = $ws.Cells.Item(30, 4)
$d.NumberFormat = "@"
$d.Value = "0.0$([char]0x2083)0862"
$d.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -9.32%  "

$ws.Cells.Item(31, 5).Value = "  -0.01%  "

$d = $ws.Cells.Item(32, 4)
$d.NumberFormat = "@"
$d.Value = "1.67"
$d.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -3.18%  "

$ws.Cells.Item(33, 5).Value = "  -4.60%  "

$d = $ws.Cells.Item(34, 4)
$d.NumberFormat = "@"
$d.Value = "148.17"
$d.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -7.04%  "

$ws.Cells.Item(35, 5).Value = "  -6.72%  "

$ws.Cells.Item(36, 5).Value = "  -7.28%  "

$ws.Cells.Item(37, 5).Value = "  -7.25%  "

$ws.Cells.Item(38, 5).Value = "  -6.47%  "

$d = $ws.Cells.Item(39, 4)
$d.NumberFormat = "@"
$d.Value = "37.97"
$d.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +1.98%  "

$ws.Cells.Item(40, 5).Value = "  -4.52%  "

$d = $ws.Cells.Item(41, 4)
$d.NumberFormat = "@"
$d.Value = "2.330.11"
$d.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -5.02%  "

$d = $ws.Cells.Item(42, 4)
$d.NumberFormat = "@"
$d.Value = "3.67"
$d.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -6.51%  "

$d = $ws.Cells.Item(43, 4)
$d.NumberFormat = "@"
$d.Value = "0.644"
$d.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -2.99%  "

$d = $ws.Cells.Item(44, 4)
$d.NumberFormat = "@"
$d.Value = "20.76"
$d.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -7.48%  "

$d = $ws.Cells.Item(45, 4)
$d.NumberFormat = "@"
$d.Value = "0.0574"
$d.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -3.48%  "

$ws.Cells.Item(46, 5).Value = "  +0.01%  "

$ws.Cells.Item(47, 5).Value = "  +1.64%  "

$ws.Cells.Item(48, 5).Value = "  -5.08%  "

$ws.Cells.Item(49, 5).Value = "  -1.52%  "

$ws.Cells.Item(50, 5).Value = "  -1.32%  "

$d = $ws.Cells.Item(51, 4)
$d.NumberFormat = "@"
$d.Value = "250.82"
$d.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -5.85%  "
